# This workbook's data rows (2-24) are cyclically rotated down by 3
# positions: the record that used to live at row N now lives at row
# ((N - 2 + 3) mod 23) + 2.  Equivalently, the NEW content of row R is the
# OLD content that used to sit at row ((R - 5) mod 23) + 2.  Row 1 (header)
# and row 25 (last record) are untouched.
#
# Because this is a rotation within the very same 2..24 block, the old
# values have to be staged out of the way first (otherwise the source for
# a later row would already have been overwritten by an earlier step).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 24
$rowCount = $lastDataRow - $firstDataRow + 1   # 23
$lastCol = 51                                  # column AY
$stageOffset = 1000                            # scratch rows, far away from real data

# 1) Stage a verbatim copy of every used cell in rows 2..24 into rows
#    1002..1024 so the originals survive the upcoming clear/rewrite.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
  for ($col = 1; $col -le $lastCol; $col++) {
    $src = $ws.Cells.Item($r, $col)
    $dst = $ws.Cells.Item($r + $stageOffset, $col)
    $src.Copy($dst)
  }
}

# 2) Wipe rows 2..24 completely (removes every cell, not just its value),
#    so columns absent from the eventual source row stay absent.
$ws.Range("A$firstDataRow`:AY$lastDataRow").ClearContents() | Out-Null

# 3) Write each destination row back from its rotated staged source row.
for ($newRow = $firstDataRow; $newRow -le $lastDataRow; $newRow++) {
  $oldRow = (($newRow - $firstDataRow - 3) % $rowCount + $rowCount) % $rowCount + $firstDataRow
  $stagedRow = $oldRow + $stageOffset
  for ($col = 1; $col -le $lastCol; $col++) {
    $src = $ws.Cells.Item($stagedRow, $col)
    $dst = $ws.Cells.Item($newRow, $col)
    $src.Copy($dst)
  }
}

# 4) Clean up the scratch rows used for staging.
$ws.Range("A$($firstDataRow + $stageOffset):AY$($lastDataRow + $stageOffset)").ClearContents() | Out-Null
